$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows (in descending order so row indices for subsequent deletes remain valid)
$rowsToDelete = @(23,22,21,17,15,14,11,10,9,8,6,5)
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}

$ws.Range("A1").Select()
